$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.232075572013855
$ws.Range("B1").Value = 2.318470001220703
$ws.Range("C1").Value = 3.275179862976074
$ws.Range("D1").Value = 2.102959394454956
$ws.Range("E1").Value = 1.360849142074585
